$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Fill in SDG column (I) values that were previously "-" placeholders
$ws1.Range("I2").Value = 1.2
$ws1.Range("I3").Value = 2
$ws1.Range("I4").Value = 1
$ws1.Range("I5").Value = 3.4
$ws1.Range("I6").Value = 4
$ws1.Range("I7").Value = 5
$ws1.Range("I8").Value = 4
$ws1.Range("I9").Value = 5
$ws1.Range("I10").Value = 5

# Row 9: new source and corrected DAC5/CRS value
$ws1.Range("F9").Value = "OCHA Indicator Registry, FAO"
$ws1.Range("H9").Value = 720.998

# Row 10: new source and corrected DAC5/CRS value
$ws1.Range("F10").Value = "OCHA Indicator Registry, Capacity4Dev"
$ws1.Range("H10").Value = 720.112

# Update selections to match final cursor position left by the edits
[void]$ws2.Range("A1").Select()
[void]$ws1.Range("F9").Select()
